$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Create new row 36 content first (currently empty); copy style from A5 (which has the standard header-row style)
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A36").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# 2) Shift existing rows 4-35 down to 5-36 (processed bottom-up) with literal values
$ws.Range("A36").Value = "ibes_1|ni|rnn_double｜all"
$ws.Range("B36").Value = [double]"0.009075753729992564"
$ws.Range("C36").Value = [double]"0.008687458498950898"
$ws.Range("D36").Value = [double]"0.0001756371360362956"
$ws.Range("E36").Value = [double]"0.0001906209278749646"
$ws.Range("F36").Value = [double]"0.1949495006477631"
$ws.Range("G36").Value = [double]"0.1262697819155144"
$ws.Range("H36").Value = [double]"0.2583620952127686"
$ws.Range("I36").Value = [double]"14166"

$ws.Range("A35").Value = "ibes_1|fwdepsqcut|rnn_eps｜all"
$ws.Range("B35").Value = [double]"0.009873209207006982"
$ws.Range("C35").Value = [double]"0.008687458498950898"
$ws.Range("D35").Value = [double]"0.0002181271559853543"
$ws.Range("E35").Value = [double]"0.0001906209278749646"
$ws.Range("F35").Value = [double]"0.0001922155460098995"
$ws.Range("G35").Value = [double]"0.1262697819155144"
$ws.Range("H35").Value = [double]"0.2583620952127686"
$ws.Range("I35").Value = [double]"14166"

$ws.Range("A34").Value = "ibes_1|fwdepsqcut-sector_code|dense2｜new with indi code -fix space"
$ws.Range("B34").Value = [double]"0.01021710973531049"
$ws.Range("C34").Value = [double]"0.008687458498950898"
$ws.Range("D34").Value = [double]"0.0002201065094334888"
$ws.Range("E34").Value = [double]"0.0001906209278749646"
$ws.Range("F34").Value = [double]"-0.00888035030069112"
$ws.Range("G34").Value = [double]"0.1262697819155144"
$ws.Range("H34").Value = [double]"0.2583620952127686"
$ws.Range("I34").Value = [double]"14166"

$ws.Range("A33").Value = "ibes_1|fwdepsqcut-industry_code|dense2｜new with indi code -fix space"
$ws.Range("B33").Value = [double]"0.009548220646664947"
$ws.Range("C33").Value = [double]"0.008687458498950898"
$ws.Range("D33").Value = [double]"0.0001933459958501393"
$ws.Range("E33").Value = [double]"0.0001906209278749646"
$ws.Range("F33").Value = [double]"0.113779158442072"
$ws.Range("G33").Value = [double]"0.1262697819155144"
$ws.Range("H33").Value = [double]"0.2583620952127686"
$ws.Range("I33").Value = [double]"14166"

$ws.Range("A32").Value = "ibes_1|fwdepsqcut|dense2｜new with indi code -fix space"
$ws.Range("B32").Value = [double]"0.009225854121736928"
$ws.Range("C32").Value = [double]"0.008687458498950898"
$ws.Range("D32").Value = [double]"0.0001822661374906326"
$ws.Range("E32").Value = [double]"0.0001906209278749646"
$ws.Range("F32").Value = [double]"0.1645648049538074"
$ws.Range("G32").Value = [double]"0.1262697819155144"
$ws.Range("H32").Value = [double]"0.2583620952127686"
$ws.Range("I32").Value = [double]"14166"

$ws.Range("A31").Value = "ibes_2|fwdepsqcut|ibes_industry -sp500"
$ws.Range("B31").Value = [double]"0.005697104856796651"
$ws.Range("C31").Value = [double]"0.005645917892372011"
$ws.Range("D31").Value = [double]"8.879743897162471E-05"
$ws.Range("E31").Value = [double]"9.170236601052772E-05"
$ws.Range("F31").Value = [double]"0.283710677227887"
$ws.Range("G31").Value = [double]"0.260277926852472"
$ws.Range("H31").Value = [double]"0.4208397194991282"
$ws.Range("I31").Value = [double]"6771"

$ws.Range("A30").Value = "ibes_2|fwdepsqcut|dense2｜new industry model -fix space"
$ws.Range("B30").Value = [double]"0.009513742997451299"
$ws.Range("C30").Value = [double]"0.008635508151540194"
$ws.Range("D30").Value = [double]"0.0002127091144161235"
$ws.Range("E30").Value = [double]"0.0001921541471413092"
$ws.Range("F30").Value = [double]"0.06095386574396355"
$ws.Range("G30").Value = [double]"0.15169780312604"
$ws.Range("H30").Value = [double]"0.2583620952127688"
$ws.Range("I30").Value = [double]"14166"

$ws.Range("A29").Value = "ibes_1|ni|cnn_rnn｜small_training_False_0"
$ws.Range("B29").Value = [double]"0.009022455578369986"
$ws.Range("C29").Value = [double]"0.008687458498950898"
$ws.Range("D29").Value = [double]"0.0001738441260953471"
$ws.Range("E29").Value = [double]"0.0001906209278749646"
$ws.Range("F29").Value = [double]"0.2031679422648356"
$ws.Range("G29").Value = [double]"0.1262697819155144"
$ws.Range("H29").Value = [double]"0.2583620952127686"
$ws.Range("I29").Value = [double]"14166"

$ws.Range("A28").Value = "ibes_2|ni|ibes_new industry_all x -indi space"
$ws.Range("B28").Value = [double]"0.007587864829138694"
$ws.Range("C28").Value = [double]"0.008635508151540194"
$ws.Range("D28").Value = [double]"0.0001393500412889987"
$ws.Range("E28").Value = [double]"0.0001921541471413092"
$ws.Range("F28").Value = [double]"0.3848118923345282"
$ws.Range("G28").Value = [double]"0.15169780312604"
$ws.Range("H28").Value = [double]"0.2583620952127688"
$ws.Range("I28").Value = [double]"14166"

$ws.Range("A27").Value = "ibes_1|fwdepsqcut-sector_code|ibes_entire_only ws -smaller space"
$ws.Range("B27").Value = [double]"0.008418017610370062"
$ws.Range("C27").Value = [double]"0.008687458498950898"
$ws.Range("D27").Value = [double]"0.0001595007627511497"
$ws.Range("E27").Value = [double]"0.0001906209278749646"
$ws.Range("F27").Value = [double]"0.2689121925027265"
$ws.Range("G27").Value = [double]"0.1262697819155144"
$ws.Range("H27").Value = [double]"0.2583620952127686"
$ws.Range("I27").Value = [double]"14166"

$ws.Range("A26").Value = "ibes_1|fwdepsqcut-industry_code|ibes_entire_only ws -smaller space"
$ws.Range("B26").Value = [double]"0.008443125284138316"
$ws.Range("C26").Value = [double]"0.008687458498950898"
$ws.Range("D26").Value = [double]"0.0001601555041560421"
$ws.Range("E26").Value = [double]"0.0001906209278749646"
$ws.Range("F26").Value = [double]"0.2659111193421716"
$ws.Range("G26").Value = [double]"0.1262697819155144"
$ws.Range("H26").Value = [double]"0.2583620952127686"
$ws.Range("I26").Value = [double]"14166"

$ws.Range("A25").Value = "ibes_1|fwdepsqcut|ibes_entire_only ws -smaller space"
$ws.Range("B25").Value = [double]"0.00844284418539723"
$ws.Range("C25").Value = [double]"0.008687458498950898"
$ws.Range("D25").Value = [double]"0.0001604772912864719"
$ws.Range("E25").Value = [double]"0.0001906209278749646"
$ws.Range("F25").Value = [double]"0.2644361756264868"
$ws.Range("G25").Value = [double]"0.1262697819155144"
$ws.Range("H25").Value = [double]"0.2583620952127686"
$ws.Range("I25").Value = [double]"14166"

$ws.Range("A24").Value = "ibes_2|ni|ibes_new industry_all x -mse"
$ws.Range("B24").Value = [double]"0.008305626259430711"
$ws.Range("C24").Value = [double]"0.008635508151540194"
$ws.Range("D24").Value = [double]"0.0001581971829802284"
$ws.Range("E24").Value = [double]"0.0001921541471413092"
$ws.Range("F24").Value = [double]"0.3016074861881056"
$ws.Range("G24").Value = [double]"0.15169780312604"
$ws.Range("H24").Value = [double]"0.2583620952127688"
$ws.Range("I24").Value = [double]"14166"

$ws.Range("A23").Value = "ibes_6|fwdepsqcut|ibes_sector_only ws -indi space"
$ws.Range("B23").Value = [double]"0.008271263221934944"
$ws.Range("C23").Value = [double]"0.008659154831887702"
$ws.Range("D23").Value = [double]"0.000162539775585851"
$ws.Range("E23").Value = [double]"0.0001949633210112301"
$ws.Range("F23").Value = [double]"0.2920577529107484"
$ws.Range("G23").Value = [double]"0.1508369500376733"
$ws.Range("H23").Value = [double]"0.2583620952127688"
$ws.Range("I23").Value = [double]"14166"

$ws.Range("A22").Value = "ibes_1|fwdepsqcut|dense2｜fix_space -best_col 10 -code 0"
$ws.Range("B22").Value = [double]"0.009537818792305376"
$ws.Range("C22").Value = [double]"0.008331870765561454"
$ws.Range("D22").Value = [double]"0.0001882493670466495"
$ws.Range("E22").Value = [double]"0.0001864150531956526"
$ws.Range("F22").Value = [double]"0.1478472391608965"
$ws.Range("G22").Value = [double]"0.1561506700669093"
$ws.Range("H22").Value = [double]"0.4578656362182956"
$ws.Range("I22").Value = [double]"2594"

$ws.Range("A21").Value = "ibes_1|fwdepsqcut-sector_code|dense2｜top15 -small space"
$ws.Range("B21").Value = [double]"0.009511178178242779"
$ws.Range("C21").Value = [double]"0.008679017950585495"
$ws.Range("D21").Value = [double]"0.0001908835104177435"
$ws.Range("E21").Value = [double]"0.0001903704215854449"
$ws.Range("F21").Value = [double]"0.1242610773958542"
$ws.Range("G21").Value = [double]"0.1266150358923953"
$ws.Range("H21").Value = [double]"0.258876568546933"
$ws.Range("I21").Value = [double]"14156"

$ws.Range("A20").Value = "ibes_1|fwdepsqcut|dense2｜top15 -small space"
$ws.Range("B20").Value = [double]"0.009500780046706191"
$ws.Range("C20").Value = [double]"0.008679017950585495"
$ws.Range("D20").Value = [double]"0.0001906222323622052"
$ws.Range("E20").Value = [double]"0.0001903704215854449"
$ws.Range("F20").Value = [double]"0.1254597737230357"
$ws.Range("G20").Value = [double]"0.1266150358923953"
$ws.Range("H20").Value = [double]"0.258876568546933"
$ws.Range("I20").Value = [double]"14156"

$ws.Range("A19").Value = "ibes_2|fwdepsqcut|ibes_new industry_only ws -indi space3 (compare using old)"
$ws.Range("B19").Value = [double]"0.00698011858859047"
$ws.Range("C19").Value = [double]"0.01116184410920933"
$ws.Range("D19").Value = [double]"0.0001206003995053095"
$ws.Range("E19").Value = [double]"0.0002568809119206026"
$ws.Range("F19").Value = [double]"0.1961572221345995"
$ws.Range("G19").Value = [double]"-0.7121988539495976"
$ws.Range("H19").Value = [double]"-0.0971244535465241"
$ws.Range("I19").Value = [double]"1545"

$ws.Range("A18").Value = "ibes_1|fwdepsqcut-sector_code|dense2｜new with indi code -fix space_sp500"
$ws.Range("B18").Value = [double]"0.007232949133456665"
$ws.Range("C18").Value = [double]"0.005640138541351778"
$ws.Range("D18").Value = [double]"0.0001213274276709267"
$ws.Range("E18").Value = [double]"9.104237776586397E-05"
$ws.Range("F18").Value = [double]"-0.01962258904981384"
$ws.Range("G18").Value = [double]"0.2348896971370824"
$ws.Range("H18").Value = [double]"0.4208397194991282"
$ws.Range("I18").Value = [double]"6771"

$ws.Range("A17").Value = "ibes_1|fwdepsqcut-industry_code|dense2｜new with indi code -fix space_sp500"
$ws.Range("B17").Value = [double]"0.006748081351245411"
$ws.Range("C17").Value = [double]"0.005640138541351778"
$ws.Range("D17").Value = [double]"0.0001082335137856306"
$ws.Range("E17").Value = [double]"9.104237776586397E-05"
$ws.Range("F17").Value = [double]"0.09041724804895079"
$ws.Range("G17").Value = [double]"0.2348896971370824"
$ws.Range("H17").Value = [double]"0.4208397194991282"
$ws.Range("I17").Value = [double]"6771"

$ws.Range("A16").Value = "ibes_1|fwdepsqcut|dense2｜new with indi code -fix space_sp500"
$ws.Range("B16").Value = [double]"0.006498537234416675"
$ws.Range("C16").Value = [double]"0.005640138541351778"
$ws.Range("D16").Value = [double]"0.0001029319227940973"
$ws.Range("E16").Value = [double]"9.104237776586397E-05"
$ws.Range("F16").Value = [double]"0.1349712457446069"
$ws.Range("G16").Value = [double]"0.2348896971370824"
$ws.Range("H16").Value = [double]"0.4208397194991282"
$ws.Range("I16").Value = [double]"6771"

$ws.Range("A15").Value = "ibes_2|fwdepsqcut|xgb xgb_space -sample_type industry -x_type fwdepsqcut"
$ws.Range("B15").Value = [double]"0.009187624256106797"
$ws.Range("C15").Value = [double]"0.007760626981832584"
$ws.Range("D15").Value = [double]"0.000170123123921569"
$ws.Range("E15").Value = [double]"0.0001663066979293028"
$ws.Range("F15").Value = [double]"0.1811552644422214"
$ws.Range("G15").Value = [double]"0.1995246680857482"
$ws.Range("H15").Value = [double]"0.24728778214858"
$ws.Range("I15").Value = [double]"140"

$ws.Range("A14").Value = "ibes_2|fwdepsqcut|ibes_industry_all x -exclude_stock"
$ws.Range("B14").Value = [double]"0.008657117254708444"
$ws.Range("C14").Value = [double]"0.008635508151540194"
$ws.Range("D14").Value = [double]"0.0001731025162188118"
$ws.Range("E14").Value = [double]"0.0001921541471413092"
$ws.Range("F14").Value = [double]"0.2358049671192343"
$ws.Range("G14").Value = [double]"0.15169780312604"
$ws.Range("H14").Value = [double]"0.2583620952127688"
$ws.Range("I14").Value = [double]"14166"

$ws.Range("A13").Value = "ibes_2|fwdepsqcut|ibes_new industry_only ws -indi space3_sp500"
$ws.Range("B13").Value = [double]"0.00566147729919014"
$ws.Range("C13").Value = [double]"0.005645917892372011"
$ws.Range("D13").Value = [double]"8.617292926780686E-05"
$ws.Range("E13").Value = [double]"9.170236601052772E-05"
$ws.Range("F13").Value = [double]"0.3048814260707366"
$ws.Range("G13").Value = [double]"0.260277926852472"
$ws.Range("H13").Value = [double]"0.4208397194991282"
$ws.Range("I13").Value = [double]"6771"

$ws.Range("A12").Value = "ibes_1|fwdepsqcut|dense2｜small_space -best_col 10 -code 0"
$ws.Range("B12").Value = [double]"0.009796859060637545"
$ws.Range("C12").Value = [double]"0.00916128901404265"
$ws.Range("D12").Value = [double]"0.0001979826847859096"
$ws.Range("E12").Value = [double]"0.0002092462997132058"
$ws.Range("F12").Value = [double]"0.07836738039900737"
$ws.Range("G12").Value = [double]"0.02593393177269743"
$ws.Range("H12").Value = [double]"0.1106454438525972"
$ws.Range("I12").Value = [double]"4659"

$ws.Range("A11").Value = "ibes_1|fwdepsqcut|dense2｜all x 0 -fix space"
$ws.Range("B11").Value = [double]"0.009081839591891672"
$ws.Range("C11").Value = [double]"0.008687458498950898"
$ws.Range("D11").Value = [double]"0.000180385026887178"
$ws.Range("E11").Value = [double]"0.0001906209278749646"
$ws.Range("F11").Value = [double]"0.1731870648290482"
$ws.Range("G11").Value = [double]"0.1262697819155144"
$ws.Range("H11").Value = [double]"0.2583620952127686"
$ws.Range("I11").Value = [double]"14166"

$ws.Range("A10").Value = "ibes_1|ni-sector_code|cnn_rnn｜all"
$ws.Range("B10").Value = [double]"0.008991444386632417"
$ws.Range("C10").Value = [double]"0.008687458498950898"
$ws.Range("D10").Value = [double]"0.0001734042650277677"
$ws.Range("E10").Value = [double]"0.0001906209278749646"
$ws.Range("F10").Value = [double]"0.2051840897612696"
$ws.Range("G10").Value = [double]"0.1262697819155144"
$ws.Range("H10").Value = [double]"0.2583620952127686"
$ws.Range("I10").Value = [double]"14166"

$ws.Range("A9").Value = "ibes_1|fwdepsqcut|cnn_rnn｜all"
$ws.Range("B9").Value = [double]"0.009193166387152663"
$ws.Range("C9").Value = [double]"0.008687458498950898"
$ws.Range("D9").Value = [double]"0.0001790877823978381"
$ws.Range("E9").Value = [double]"0.0001906209278749646"
$ws.Range("F9").Value = [double]"0.1791331155760224"
$ws.Range("G9").Value = [double]"0.1262697819155144"
$ws.Range("H9").Value = [double]"0.2583620952127686"
$ws.Range("I9").Value = [double]"14166"

$ws.Range("A8").Value = "ibes_1|fwdepsqcut-sector_code|cnn_rnn｜adj_space_True_2"
$ws.Range("B8").Value = [double]"0.009050797912742738"
$ws.Range("C8").Value = [double]"0.008687458498950898"
$ws.Range("D8").Value = [double]"0.0001747927346526816"
$ws.Range("E8").Value = [double]"0.0001906209278749646"
$ws.Range("F8").Value = [double]"0.1988198994191919"
$ws.Range("G8").Value = [double]"0.1262697819155144"
$ws.Range("H8").Value = [double]"0.2583620952127686"
$ws.Range("I8").Value = [double]"14166"

$ws.Range("A7").Value = "ibes_2|fwdepsqcut|ibes_new industry_only ws -indi space3"
$ws.Range("B7").Value = [double]"0.008303107575837335"
$ws.Range("C7").Value = [double]"0.008635508151540194"
$ws.Range("D7").Value = [double]"0.0001612992868167532"
$ws.Range("E7").Value = [double]"0.0001921541471413092"
$ws.Range("F7").Value = [double]"0.2879126399482278"
$ws.Range("G7").Value = [double]"0.15169780312604"
$ws.Range("H7").Value = [double]"0.2583620952127688"
$ws.Range("I7").Value = [double]"14166"

$ws.Range("A6").Value = "ibes_6|ni|ibes_sector_all x"
$ws.Range("B6").Value = [double]"0.007727376897992178"
$ws.Range("C6").Value = [double]"0.008659154831887702"
$ws.Range("D6").Value = [double]"0.0001456055809025543"
$ws.Range("E6").Value = [double]"0.0001949633210112301"
$ws.Range("F6").Value = [double]"0.3658146643715235"
$ws.Range("G6").Value = [double]"0.1508369500376733"
$ws.Range("H6").Value = [double]"0.2583620952127688"
$ws.Range("I6").Value = [double]"14166"

$ws.Range("A5").Value = "ibes_1|ni-sector_code|rnn_double｜adj_space__exclude_fwd_2"
$ws.Range("B5").Value = [double]"0.008996341110758196"
$ws.Range("C5").Value = [double]"0.008687458498950898"
$ws.Range("D5").Value = [double]"0.0001736821065404229"
$ws.Range("E5").Value = [double]"0.0001906209278749646"
$ws.Range("F5").Value = [double]"0.2039105752098954"
$ws.Range("G5").Value = [double]"0.1262697819155144"
$ws.Range("H5").Value = [double]"0.2583620952127686"
$ws.Range("I5").Value = [double]"14166"

# 3) Write the brand-new row 4 data
$ws.Range("A4").Value = "ibes_1|fwdepsqcut|dense2｜small_space -best_col 15 -code 0 -exclude_fwd True"
$ws.Range("B4").Value = [double]"0.009913389281997772"
$ws.Range("C4").Value = [double]"0.009030432768250784"
$ws.Range("D4").Value = [double]"0.000197928676326859"
$ws.Range("E4").Value = [double]"0.0002181746930107264"
$ws.Range("F4").Value = [double]"0.1982174071110625"
$ws.Range("G4").Value = [double]"0.116203501628984"
$ws.Range("H4").Value = [double]"0.4935486705331952"
$ws.Range("I4").Value = [double]"629"
